$wb = $excel.ActiveWorkbook

$item = $wb.Worksheets.Item("Item")
$item.Range("X5").Value = "72h3m0.5s"
$item.Range("X3").Value = "1m0.5s"
$item.Range("X4").Value = "1h1s"

$activity = $wb.Worksheets.Item("Activity")
$activity.Range("L4").Value = "72h3m0.5s"
$activity.Range("L2").Value = "1m0.5s"
$activity.Range("L3").Value = "1h1s"
$activity.Range("L6").Value = "1h1m"
$activity.Range("L5").Value = "1h1m1s"

$env = $wb.Worksheets.Item("Env")
$env.Range("A1").Value = "ZoneID"
$env.Range("A2").Value = "WorldID"

$activity.PageSetup.PaperSize = 9
$activity.PageSetup.Orientation = 1

[void]$item.Range("X3:X5").Select()
[void]$activity.Activate()
[void]$activity.Range("L9").Select()
